# Auto-generated edit script for data-dictionary.xlsx update
$wb = $excel.ActiveWorkbook

# --- Sheet 'customers': update loc_lead_source / loc_lead_source_type descriptions ---
$wsCustomers = $wb.Worksheets.Item("customers")
$wsCustomers.Range("C10").Value = 'marketing lead source for customer (level 1 depth, ex. "radio")'
$wsCustomers.Range("C11").Value = 'marketing lead source customer (level 2 depth, ex. "radio for discounted tune-up in 2017")'
[void]$wsCustomers.Range("C22").Select()

# --- Sheet 'sales_calls': populate data dictionary table ---
$wsSalesCalls = $wb.Worksheets.Item("sales_calls")
$sheet2Data = @(
  @('column_name', 'data_type', 'description'),
  @('billing_acct_id', 'int', 'unique account identifier; primary key'),
  @('call_type', 'varchar', 'category of estimate (furnace, water heater, etc.)'),
  @('location_id', 'int', 'unique location identifier'),
  @('taken_date', 'datetime', 'date call was scheduled'),
  @('job_lead_source', 'varchar', 'marketing lead source for particular sales call (level 1 depth, ex. "technician")'),
  @('job_lead_source_type', 'varchar', 'marketing lead source particular sales call (level 2 depth, ex. "technician- Jairo")'),
  @('salesperson', 'varchar', 'salesperson'),
  @('call_date', 'datetime', 'date of appointment'),
  @('call_time', 'timestamps', 'time of appointment'),
)
for ($r = 0; $r -lt $sheet2Data.Length; $r++) {
  for ($c = 0; $c -lt $sheet2Data[$r].Length; $c++) {
    $wsSalesCalls.Cells.Item($r + 1, $c + 1).Value = $sheet2Data[$r][$c]
  }
}
$wsSalesCalls.Columns.Item(1).ColumnWidth = 18.334
$wsSalesCalls.Columns.Item(2).ColumnWidth = 14.668
$wsSalesCalls.Columns.Item(3).ColumnWidth = 97.001
[void]$wsSalesCalls.Range("B4:C4").Select()

# --- Sheet 'job_summary': populate data dictionary table ---
$wsJobSummary = $wb.Worksheets.Item("job_summary")
$sheet3Data = @(
  @('column_name', 'data_type', 'description'),
  @('billing_acct_id', 'int', 'unique account identifier; primary key'),
  @('job_no', 'int', 'unique job identifier; secondary key'),
  @('end_date', 'datetime', 'date job ended'),
  @('location_id', 'int', 'unique location identifier'),
  @('job_class', 'varchar', 'class of job (water heater install, etc.)'),
  @('job_type', 'varchar', 'type of job (water heater install, etc.)'),
  @('dept', 'dept', 'department- install, service, maintenance, etc.'),
)
for ($r = 0; $r -lt $sheet3Data.Length; $r++) {
  for ($c = 0; $c -lt $sheet3Data[$r].Length; $c++) {
    $wsJobSummary.Cells.Item($r + 1, $c + 1).Value = $sheet3Data[$r][$c]
  }
}
$wsJobSummary.Columns.Item(1).ColumnWidth = 20.501
$wsJobSummary.Columns.Item(2).ColumnWidth = 15.835
$wsJobSummary.Columns.Item(3).ColumnWidth = 82.334
[void]$wsJobSummary.Range("C9").Select()

[void]$wsCustomers.Select()
